$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: insert a new leading blank paragraph before the document's
# very first paragraph ("Given an array ...").
#
# We do this via Find/Replace (wdReplaceOne) rather than
# Range(0,0).InsertParagraphBefore() because, at the absolute start of
# the story, InsertParagraphBefore() on this engine ends up inheriting
# the Bold run-formatting that appears later in that first paragraph
# (the bolded "arr[]"), leaving a stray <w:rPr><w:b/></w:rPr> on the
# new empty paragraph. Prefixing the matched text with a paragraph
# mark ("^p") via Find/Replace instead yields a clean empty paragraph.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Given an array", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "^pGiven an array", 1) | Out-Null

# ------------------------------------------------------------------
# Step 2: the document used to contain the whole problem statement
# block (statement + two examples, no explanations) twice in a row;
# drop the first (now second-of-original, i.e. paragraphs 2..9 after
# the insert above) copy, keeping the second copy -- the one that has
# the "Explanation:" paragraphs -- in place.
# ------------------------------------------------------------------
$dupStart = $d.Paragraphs(2).Range.Start
$dupEnd = $d.Paragraphs(9).Range.End
$d.Range($dupStart, $dupEnd).Delete() | Out-Null

# ------------------------------------------------------------------
# Step 3: append the Java solution listing (plus two blank separator
# paragraphs) right after the final "Explanation: ..." paragraph and
# before the document's trailing blank paragraph, so that trailing
# blank paragraph stays the very last paragraph in the document.
# ------------------------------------------------------------------

# 3a: two clean blank paragraphs, built with InsertParagraphBefore (on
# this engine this yields a bare <w:p><w:r/></w:p> here, unlike a
# "\r" sent through InsertBefore text which leaves an explicit empty
# <w:t></w:t> behind).
$lastIndex = $d.Paragraphs.Count
$insertPos = $d.Paragraphs($lastIndex).Range.Start
$d.Range($insertPos, $insertPos).InsertParagraphBefore() | Out-Null
$insertPos = $d.Paragraphs($d.Paragraphs.Count).Range.Start
$d.Range($insertPos, $insertPos).InsertParagraphBefore() | Out-Null

# 3b: the Java code listing itself, typed as one multi-line block so
# each "`r" becomes its own paragraph.
$insertPos = $d.Paragraphs($d.Paragraphs.Count).Range.Start
$codeBlock = "class Solution {`r" + `
    "    public int maximumDistance(int [] nums , int n) {`r" + `
    "        Map<Integer , Integer> map = new HashMap<>();`r" + `
    "        int max = -1;`r" + `
    "        for(int  i = 0 ;  i< n ; i++){`r" + `
    "            if(map.contains(nums[i])){`r" + `
    "                max = math.max(max , i - map.get(nums[i]));`r" + `
    "            }else{`r" + `
    "                map.put(nums[i] , i);`r" + `
    "            }`r" + `
    "        }`r" + `
    "        `r" + `
    "    }`r" + `
    "}`r"
$d.Range($insertPos, $insertPos).InsertBefore($codeBlock) | Out-Null

Write-Output "Final paragraph count:" $d.Paragraphs.Count
